# Weekly refresh of the "Espinaca" (spinach) price series:
# a new daily observation is inserted as row 550 (pushing the existing
# rows 550-592 down to 551-593), extending the used range to A1:R593.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 550; Excel shifts rows 550:592 down to 551:593.
$ws.Rows("550").Insert()

# Populate the newly inserted row with the new market observation.
$ws.Range("A550").Value = 6
$ws.Range("B550").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C550").Value = "Metropolitana"
$ws.Range("D550").Value = 44783
$ws.Range("D550").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E550").Value = 13
$ws.Range("F550").Value = 100112012
$ws.Range("G550").Value = "Espinaca"
$ws.Range("H550").Value = "Sin especificar"
$ws.Range("I550").Value = "Primera"
$ws.Range("J550").Value = 480
$ws.Range("K550").Value = 6000
$ws.Range("L550").Value = 6500
$ws.Range("M550").Value = 6219
$ws.Range("N550").Value = "$/cuna 10 kilos"
$ws.Range("O550").Value = "Provincia de Chacabuco"
$ws.Range("P550").Value = 622
$ws.Range("Q550").Value = 10
$ws.Range("R550").Value = "Hortaliza"
